$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes current rows 6-8 down to 7-9),
# carrying formatting along like Excel's native row insert.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with this week's data point.
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44767
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 100112036
$ws.Range("G6").Value = "Caigua"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 20000
$ws.Range("N6").Value = '$/caja 15 kilos'
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 1333
$ws.Range("Q6").Value = 15
$ws.Range("R6").Value = "Hortaliza"
